$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.1
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2.05
$ws.Range("X2").Value = 1.18
